$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 5.8
$ws.Cells.Item(2, 8).Value = 1.77
$ws.Cells.Item(2, 10).Value = 3.5
$ws.Cells.Item(2, 11).Value = 4
$ws.Cells.Item(2, 12).Value = 1.43
$ws.Cells.Item(2, 14).Value = 3.05
$ws.Cells.Item(2, 15).Value = 1.38
$ws.Cells.Item(2, 16).Value = 1.72
$ws.Cells.Item(2, 17).Value = 2.1
$ws.Cells.Item(2, 18).Value = 1.27
$ws.Cells.Item(2, 19).Value = 3.85
$ws.Cells.Item(2, 20).Value = 2.02
$ws.Cells.Item(2, 21).Value = 1.8
$ws.Cells.Item(2, 23).Value = 1.21
$ws.Cells.Item(2, 24).Value = 23
$ws.Cells.Item(2, 25).Value = 14
$ws.Cells.Item(2, 27).Value = 900
$ws.Cells.Item(2, 28).Value = 1000
$ws.Cells.Item(2, 29).Value = 14
$ws.Cells.Item(2, 30).Value = 21
$ws.Cells.Item(2, 32).Value = 980
$ws.Cells.Item(2, 33).Value = 980
$ws.Cells.Item(2, 36).Value = 900
$ws.Cells.Item(2, 39).Value = 580
$ws.Cells.Item(2, 40).Value = 600

# Row 3
$ws.Cells.Item(3, 6).Value = 1.85
$ws.Cells.Item(3, 8).Value = 4.4
$ws.Cells.Item(3, 9).Value = 5.3
$ws.Cells.Item(3, 12).Value = 1.35
$ws.Cells.Item(3, 14).Value = 3.45
$ws.Cells.Item(3, 16).Value = 1.84
$ws.Cells.Item(3, 17).Value = 1.87
$ws.Cells.Item(3, 18).Value = 1.31
$ws.Cells.Item(3, 21).Value = 2
$ws.Cells.Item(3, 22).Value = 1.23
$ws.Cells.Item(3, 23).Value = 1.99
$ws.Cells.Item(3, 26).Value = 980
$ws.Cells.Item(3, 27).Value = 900
$ws.Cells.Item(3, 31).Value = 190
$ws.Cells.Item(3, 32).Value = 25
$ws.Cells.Item(3, 34).Value = 60
$ws.Cells.Item(3, 35).Value = 500
$ws.Cells.Item(3, 36).Value = 48
$ws.Cells.Item(3, 37).Value = 48
$ws.Cells.Item(3, 38).Value = 980
$ws.Cells.Item(3, 39).Value = 580

# Row 4
$ws.Cells.Item(4, 7).Value = 1.36
$ws.Cells.Item(4, 8).Value = 7.2
$ws.Cells.Item(4, 14).Value = 1.02
$ws.Cells.Item(4, 16).Value = 3.05
$ws.Cells.Item(4, 17).Value = 1.45
$ws.Cells.Item(4, 18).Value = 1.77
$ws.Cells.Item(4, 19).Value = 2
$ws.Cells.Item(4, 20).Value = 1.67
$ws.Cells.Item(4, 21).Value = 2.06
$ws.Cells.Item(4, 23).Value = 3.65
$ws.Cells.Item(4, 24).Value = 980
$ws.Cells.Item(4, 28).Value = 1000
$ws.Cells.Item(4, 29).Value = 1000
$ws.Cells.Item(4, 30).Value = 980
$ws.Cells.Item(4, 32).Value = 1000
$ws.Cells.Item(4, 33).Value = 1000
$ws.Cells.Item(4, 34).Value = 980
$ws.Cells.Item(4, 37).Value = 1000
$ws.Cells.Item(4, 38).Value = 980
$ws.Cells.Item(4, 40).Value = 6.2

# Row 5
$ws.Cells.Item(5, 23).Value = 1.5

# Row 6
$ws.Cells.Item(6, 6).Value = 1.84
$ws.Cells.Item(6, 7).Value = 1.94
$ws.Cells.Item(6, 10).Value = 3.75
$ws.Cells.Item(6, 12).Value = 1.32
$ws.Cells.Item(6, 22).Value = 1.26
$ws.Cells.Item(6, 23).Value = 2.06
$ws.Cells.Item(6, 27).Value = 900
$ws.Cells.Item(6, 38).Value = 130
$ws.Cells.Item(6, 39).Value = 580

# Row 7
$ws.Cells.Item(7, 6).Value = 1.9
$ws.Cells.Item(7, 9).Value = 4.7
$ws.Cells.Item(7, 14).Value = 3.45
$ws.Cells.Item(7, 15).Value = 1.3
$ws.Cells.Item(7, 16).Value = 1.9
$ws.Cells.Item(7, 20).Value = 1.79
$ws.Cells.Item(7, 22).Value = 1.29
$ws.Cells.Item(7, 25).Value = 970
$ws.Cells.Item(7, 26).Value = 970
$ws.Cells.Item(7, 28).Value = 21
$ws.Cells.Item(7, 29).Value = 16
$ws.Cells.Item(7, 30).Value = 970
$ws.Cells.Item(7, 31).Value = 500
$ws.Cells.Item(7, 32).Value = 24
$ws.Cells.Item(7, 33).Value = 22
$ws.Cells.Item(7, 34).Value = 970
$ws.Cells.Item(7, 35).Value = 500
$ws.Cells.Item(7, 36).Value = 970
$ws.Cells.Item(7, 37).Value = 970
$ws.Cells.Item(7, 38).Value = 970
$ws.Cells.Item(7, 40).Value = 55
$ws.Cells.Item(7, 41).Value = 500

# Row 8
$ws.Cells.Item(8, 10).Value = 7.6
$ws.Cells.Item(8, 11).Value = 9
$ws.Cells.Item(8, 19).Value = 2.62
$ws.Cells.Item(8, 20).Value = 2.66
$ws.Cells.Item(8, 21).Value = 1.52
$ws.Cells.Item(8, 32).Value = 12.5
$ws.Cells.Item(8, 40).Value = 15

# Row 9
$ws.Cells.Item(9, 7).Value = 1.84
$ws.Cells.Item(9, 8).Value = 4.5
$ws.Cells.Item(9, 10).Value = 3.9
$ws.Cells.Item(9, 11).Value = 4.3
$ws.Cells.Item(9, 23).Value = 2.18

# Row 10
$ws.Cells.Item(10, 6).Value = 2.36
$ws.Cells.Item(10, 9).Value = 3.75
$ws.Cells.Item(10, 25).Value = 1000
$ws.Cells.Item(10, 28).Value = 1000
$ws.Cells.Item(10, 29).Value = 14
$ws.Cells.Item(10, 30).Value = 32
$ws.Cells.Item(10, 32).Value = 34

# Row 11
$ws.Cells.Item(11, 14).Value = 2.78
$ws.Cells.Item(11, 19).Value = 5.5
$ws.Cells.Item(11, 29).Value = 7

# Row 12
$ws.Cells.Item(12, 14).Value = 4.9
$ws.Cells.Item(12, 16).Value = 2.36
$ws.Cells.Item(12, 26).Value = 65
$ws.Cells.Item(12, 27).Value = 900
$ws.Cells.Item(12, 30).Value = 15
$ws.Cells.Item(12, 31).Value = 85
$ws.Cells.Item(12, 35).Value = 95
$ws.Cells.Item(12, 36).Value = 32
$ws.Cells.Item(12, 38).Value = 29
$ws.Cells.Item(12, 39).Value = 580

# Row 13
$ws.Cells.Item(13, 6).Value = 2.18
$ws.Cells.Item(13, 7).Value = 2.32
$ws.Cells.Item(13, 8).Value = 3.1
$ws.Cells.Item(13, 9).Value = 3.35
$ws.Cells.Item(13, 11).Value = 4.2
$ws.Cells.Item(13, 16).Value = 3.05
$ws.Cells.Item(13, 19).Value = 1.98
$ws.Cells.Item(13, 20).Value = 1.44
$ws.Cells.Item(13, 21).Value = 3
$ws.Cells.Item(13, 23).Value = 1.75
$ws.Cells.Item(13, 24).Value = 100
$ws.Cells.Item(13, 30).Value = 30
$ws.Cells.Item(13, 35).Value = 75
$ws.Cells.Item(13, 36).Value = 1000
$ws.Cells.Item(13, 37).Value = 50
$ws.Cells.Item(13, 39).Value = 580
$ws.Cells.Item(13, 40).Value = 9.4

# Row 14
$ws.Cells.Item(14, 7).Value = 2.2
$ws.Cells.Item(14, 9).Value = 3.4
$ws.Cells.Item(14, 10).Value = 4.2
$ws.Cells.Item(14, 14).Value = 6.6
$ws.Cells.Item(14, 21).Value = 2.92
$ws.Cells.Item(14, 22).Value = 1.42
$ws.Cells.Item(14, 27).Value = 500
$ws.Cells.Item(14, 30).Value = 32
$ws.Cells.Item(14, 33).Value = 27
$ws.Cells.Item(14, 34).Value = 32
$ws.Cells.Item(14, 39).Value = 580

# Row 15
$ws.Cells.Item(15, 7).Value = 1.91
$ws.Cells.Item(15, 8).Value = 3.85
$ws.Cells.Item(15, 10).Value = 4.3
$ws.Cells.Item(15, 11).Value = 4.7
$ws.Cells.Item(15, 17).Value = 1.44
$ws.Cells.Item(15, 18).Value = 1.79
$ws.Cells.Item(15, 19).Value = 2.08
$ws.Cells.Item(15, 20).Value = 1.49
$ws.Cells.Item(15, 21).Value = 2.72
$ws.Cells.Item(15, 24).Value = 34
$ws.Cells.Item(15, 26).Value = 95
$ws.Cells.Item(15, 28).Value = 17.5
$ws.Cells.Item(15, 31).Value = 95
$ws.Cells.Item(15, 32).Value = 17.5
$ws.Cells.Item(15, 34).Value = 16
$ws.Cells.Item(15, 35).Value = 95
$ws.Cells.Item(15, 36).Value = 44
$ws.Cells.Item(15, 39).Value = 580

# Row 16
$ws.Cells.Item(16, 9).Value = 4.5
$ws.Cells.Item(16, 14).Value = 3.15
$ws.Cells.Item(16, 16).Value = 1.72
$ws.Cells.Item(16, 17).Value = 2.1
$ws.Cells.Item(16, 18).Value = 1.26
$ws.Cells.Item(16, 19).Value = 4.4
$ws.Cells.Item(16, 21).Value = 2
$ws.Cells.Item(16, 22).Value = 1.29
$ws.Cells.Item(16, 25).Value = 26
$ws.Cells.Item(16, 26).Value = 38
$ws.Cells.Item(16, 34).Value = 42
$ws.Cells.Item(16, 38).Value = 130
$ws.Cells.Item(16, 39).Value = 580
$ws.Cells.Item(16, 41).Value = 250

# Row 17
$ws.Cells.Item(17, 15).Value = 1.5
$ws.Cells.Item(17, 24).Value = 9
$ws.Cells.Item(17, 26).Value = 16
$ws.Cells.Item(17, 34).Value = 22
$ws.Cells.Item(17, 35).Value = 60

# Row 18
$ws.Cells.Item(18, 10).Value = 6.8
$ws.Cells.Item(18, 17).Value = 1.71
$ws.Cells.Item(18, 19).Value = 2.82
$ws.Cells.Item(18, 20).Value = 2.4
$ws.Cells.Item(18, 21).Value = 1.7
$ws.Cells.Item(18, 27).Value = 810
$ws.Cells.Item(18, 36).Value = 8.6
$ws.Cells.Item(18, 38).Value = 44
$ws.Cells.Item(18, 41).Value = 1000

# Row 19
$ws.Cells.Item(19, 9).Value = 2.84
$ws.Cells.Item(19, 12).Value = 1.46
$ws.Cells.Item(19, 16).Value = 1.81
$ws.Cells.Item(19, 17).Value = 2.2
$ws.Cells.Item(19, 18).Value = 1.31
$ws.Cells.Item(19, 19).Value = 4
$ws.Cells.Item(19, 29).Value = 7
$ws.Cells.Item(19, 35).Value = 48
$ws.Cells.Item(19, 40).Value = 32
